# Optimized searching algorithm for previous games and implemented TOCSV method.
# This script rewrites the per-champion stats table: row 2's data is replaced
# (Yorick -> Thresh) and five additional champion rows (Bard, Sion, Yorick,
# Yone, Ornn) are appended below it, each formatted to match the existing
# data row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-ChampionRow {
    param(
        [int]$Row,
        [string]$Name,
        $Kills,
        $Deaths,
        $Assists,
        $WinRate,
        $AvgCS,
        $GPM,
        [string]$KP,
        $AvgDPM,
        $AvgAbusage,
        $VisionScore,
        $NumGames
    )

    if ($Row -ne 2) {
        # Duplicate the formatting of the existing template data row (row 2)
        # onto the new row before filling in values.
        $ws.Range("A2:L2").Copy()
        $ws.Range("A" + $Row + ":L" + $Row).PasteSpecial(-4122)
    }

    $ws.Cells.Item($Row, 1).Value = $Name
    $ws.Cells.Item($Row, 2).Value = $Kills
    $ws.Cells.Item($Row, 3).Value = $Deaths
    $ws.Cells.Item($Row, 4).Value = $Assists
    $ws.Cells.Item($Row, 5).Value = $WinRate
    $ws.Cells.Item($Row, 6).Value = $AvgCS
    $ws.Cells.Item($Row, 7).Value = $GPM

    # KP is a literal text percentage (e.g. "0.0%"), not a numeric value, so
    # force text formatting before assigning it to avoid automatic percent
    # number-conversion, then drop the now-unneeded explicit format again.
    $ws.Cells.Item($Row, 8).NumberFormat = "@"
    $ws.Cells.Item($Row, 8).Value = $KP
    $ws.Cells.Item($Row, 8).ClearFormats()

    $ws.Cells.Item($Row, 9).Value = $AvgDPM
    $ws.Cells.Item($Row, 10).Value = $AvgAbusage
    $ws.Cells.Item($Row, 11).Value = $VisionScore
    $ws.Cells.Item($Row, 12).Value = $NumGames
}

Set-ChampionRow 2 "Thresh" 3 8 14 0   1 285 "0.0%" 264 89  61 2
Set-ChampionRow 3 "Bard"   0 5 3  0   1 220 "0.0%" 217 84  23 1
Set-ChampionRow 4 "Sion"   3 5 6  0   7 400 "0.0%" 640 185 14 12
Set-ChampionRow 5 "Yorick" 4 4 4  0   8 440 "0.0%" 643 216 15 52
Set-ChampionRow 6 "Yone"   4 9 8  0   7 371 "0.0%" 707 292 11 1
Set-ChampionRow 7 "Ornn"   4 2 5  100 7 366 "0.0%" 659 218 19 1
